# Update recomputed TPM-derived NATMI metrics (ligand/receptor/edge columns)
# for the Tnfsf10 -> Tnfrsf10b pair, per the new TPM values used upstream.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 25.70233233333333
$ws.Range("H2").Value = 77.106997
$ws.Range("I2").Value = 0.9446552342719222
$ws.Range("J2").Value = 0.9446552342719222
$ws.Range("M2").Value = 14.10125566666667
$ws.Range("N2").Value = 42.303767
$ws.Range("O2").Value = 0.7585903740943118
$ws.Range("P2").Value = 0.7585903740943116
$ws.Range("Q2").Value = 362.4351594619666
$ws.Range("R2").Value = 3261.916435157699
$ws.Range("S2").Value = 0.7166063675564873
$ws.Range("T2").Value = 0.7166063675564871

$ws.Range("G3").Value = 25.70233233333333
$ws.Range("H3").Value = 77.106997
$ws.Range("I3").Value = 0.9446552342719222
$ws.Range("J3").Value = 0.9446552342719222
$ws.Range("O3").Value = 0.1642962051354147
$ws.Range("P3").Value = 0.1642962051354147
$ws.Range("Q3").Value = 78.496542192937
$ws.Range("R3").Value = 706.468879736433
$ws.Range("S3").Value = 0.155203270152183
$ws.Range("T3").Value = 0.1552032701521829

$ws.Range("G4").Value = 25.70233233333333
$ws.Range("H4").Value = 77.106997
$ws.Range("I4").Value = 0.9446552342719222
$ws.Range("J4").Value = 0.9446552342719222
$ws.Range("M4").Value = 1.366842
$ws.Range("N4").Value = 4.100526
$ws.Range("O4").Value = 0.07353055703818179
$ws.Range("P4").Value = 0.07353055703818176
$ws.Range("Q4").Value = 35.131027331158
$ws.Range("R4").Value = 316.1792459804221
$ws.Range("S4").Value = 0.06946102558504856
$ws.Range("T4").Value = 0.06946102558504853

$ws.Range("G5").Value = 25.70233233333333
$ws.Range("H5").Value = 77.106997
$ws.Range("I5").Value = 0.9446552342719222
$ws.Range("J5").Value = 0.9446552342719222
$ws.Range("M5").Value = 0.06660100000000001
$ws.Range("N5").Value = 0.199803
$ws.Range("O5").Value = 0.003582863732091891
$ws.Range("P5").Value = 0.00358286373209189
$ws.Range("Q5").Value = 1.711801035732334
$ws.Range("R5").Value = 15.406209321591
$ws.Range("S5").Value = 0.003384570978203639
$ws.Range("T5").Value = 0.003384570978203638

$ws.Range("I6").Value = 0.00861995025144722
$ws.Range("J6").Value = 0.00861995025144722
$ws.Range("M6").Value = 14.10125566666667
$ws.Range("N6").Value = 42.303767
$ws.Range("O6").Value = 0.7585903740943118
$ws.Range("P6").Value = 0.7585903740943116
$ws.Range("Q6").Value = 3.307209795270333
$ws.Range("R6").Value = 29.764888157433
$ws.Range("S6").Value = 0.006539011285919705
$ws.Range("T6").Value = 0.006539011285919702

$ws.Range("I7").Value = 0.00861995025144722
$ws.Range("J7").Value = 0.00861995025144722
$ws.Range("O7").Value = 0.1642962051354147
$ws.Range("P7").Value = 0.1642962051354147
$ws.Range("S7").Value = 0.001416225114768842
$ws.Range("T7").Value = 0.001416225114768842

$ws.Range("I8").Value = 0.00861995025144722
$ws.Range("J8").Value = 0.00861995025144722
$ws.Range("M8").Value = 1.366842
$ws.Range("N8").Value = 4.100526
$ws.Range("O8").Value = 0.07353055703818179
$ws.Range("P8").Value = 0.07353055703818176
$ws.Range("Q8").Value = 0.320569554786
$ws.Range("R8").Value = 2.885125993074
$ws.Range("S8").Value = 0.0006338297436303292
$ws.Range("T8").Value = 0.000633829743630329

$ws.Range("I9").Value = 0.00861995025144722
$ws.Range("J9").Value = 0.00861995025144722
$ws.Range("M9").Value = 0.06660100000000001
$ws.Range("N9").Value = 0.199803
$ws.Range("O9").Value = 0.003582863732091891
$ws.Range("P9").Value = 0.00358286373209189
$ws.Range("Q9").Value = 0.015620132333
$ws.Range("R9").Value = 0.140581190997
$ws.Range("S9").Value = 0.00003088410712834662
$ws.Range("T9").Value = 0.00003088410712834662

$ws.Range("G10").Value = 0.2552456666666666
$ws.Range("H10").Value = 0.765737
$ws.Range("I10").Value = 0.009381216922838777
$ws.Range("J10").Value = 0.009381216922838777
$ws.Range("M10").Value = 14.10125566666667
$ws.Range("N10").Value = 42.303767
$ws.Range("O10").Value = 0.7585903740943118
$ws.Range("P10").Value = 0.7585903740943116
$ws.Range("Q10").Value = 3.599284403475444
$ws.Range("R10").Value = 32.393559631279
$ws.Range("S10").Value = 0.007116500854956157
$ws.Range("T10").Value = 0.007116500854956155

$ws.Range("G11").Value = 0.2552456666666666
$ws.Range("H11").Value = 0.765737
$ws.Range("I11").Value = 0.009381216922838777
$ws.Range("J11").Value = 0.009381216922838777
$ws.Range("O11").Value = 0.1642962051354147
$ws.Range("P11").Value = 0.1642962051354147
$ws.Range("Q11").Value = 0.7795363464769999
$ws.Range("R11").Value = 7.015827118293
$ws.Range("S11").Value = 0.001541298339974544
$ws.Range("T11").Value = 0.001541298339974543

$ws.Range("G12").Value = 0.2552456666666666
$ws.Range("H12").Value = 0.765737
$ws.Range("I12").Value = 0.009381216922838777
$ws.Range("J12").Value = 0.009381216922838777
$ws.Range("M12").Value = 1.366842
$ws.Range("N12").Value = 4.100526
$ws.Range("O12").Value = 0.07353055703818179
$ws.Range("P12").Value = 0.07353055703818176
$ws.Range("Q12").Value = 0.348880497518
$ws.Range("R12").Value = 3.139924477662
$ws.Range("S12").Value = 0.000689806106032353
$ws.Range("T12").Value = 0.0006898061060323527

$ws.Range("G13").Value = 0.2552456666666666
$ws.Range("H13").Value = 0.765737
$ws.Range("I13").Value = 0.009381216922838777
$ws.Range("J13").Value = 0.009381216922838777
$ws.Range("M13").Value = 0.06660100000000001
$ws.Range("N13").Value = 0.199803
$ws.Range("O13").Value = 0.003582863732091891
$ws.Range("P13").Value = 0.00358286373209189
$ws.Range("Q13").Value = 0.01699961664566667
$ws.Range("R13").Value = 0.152996549811
$ws.Range("S13").Value = 0.00003361162187572575
$ws.Range("T13").Value = 0.00003361162187572574

$ws.Range("G14").Value = 1.016050666666667
$ws.Range("H14").Value = 3.048152
$ws.Range("I14").Value = 0.03734359855379179
$ws.Range("J14").Value = 0.03734359855379179
$ws.Range("M14").Value = 14.10125566666667
$ws.Range("N14").Value = 42.303767
$ws.Range("O14").Value = 0.7585903740943118
$ws.Range("P14").Value = 0.7585903740943116
$ws.Range("Q14").Value = 14.32759022095378
$ws.Range("R14").Value = 128.948311988584
$ws.Range("S14").Value = 0.02832849439694872
$ws.Range("T14").Value = 0.02832849439694871

$ws.Range("G15").Value = 1.016050666666667
$ws.Range("H15").Value = 3.048152
$ws.Range("I15").Value = 0.03734359855379179
$ws.Range("J15").Value = 0.03734359855379179
$ws.Range("O15").Value = 0.1642962051354147
$ws.Range("P15").Value = 0.1642962051354147
$ws.Range("Q15").Value = 3.103082747192
$ws.Range("R15").Value = 27.927744724728
$ws.Range("S15").Value = 0.006135411528488352
$ws.Range("T15").Value = 0.006135411528488351

$ws.Range("G16").Value = 1.016050666666667
$ws.Range("H16").Value = 3.048152
$ws.Range("I16").Value = 0.03734359855379179
$ws.Range("J16").Value = 0.03734359855379179
$ws.Range("M16").Value = 1.366842
$ws.Range("N16").Value = 4.100526
$ws.Range("O16").Value = 0.07353055703818179
$ws.Range("P16").Value = 0.07353055703818176
$ws.Range("Q16").Value = 1.388780725328
$ws.Range("R16").Value = 12.499026527952
$ws.Range("S16").Value = 0.00274589560347055
$ws.Range("T16").Value = 0.002745895603470549

$ws.Range("G17").Value = 1.016050666666667
$ws.Range("H17").Value = 3.048152
$ws.Range("I17").Value = 0.03734359855379179
$ws.Range("J17").Value = 0.03734359855379179
$ws.Range("M17").Value = 0.06660100000000001
$ws.Range("N17").Value = 0.199803
$ws.Range("O17").Value = 0.003582863732091891
$ws.Range("P17").Value = 0.00358286373209189
$ws.Range("Q17").Value = 0.06766999045066667
$ws.Range("R17").Value = 0.609029914056
$ws.Range("S17").Value = 0.0001337970248841798
$ws.Range("T17").Value = 0.0001337970248841798

